$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '59.752.81'
Set-TextValue 'E2' '  +0.70%  '
Set-TextValue 'D3' '2.670.48'
Set-TextValue 'E3' '  +2.68%  '
Set-TextValue 'E4' '  -0.06%  '
Set-TextValue 'D5' '538.51'
Set-TextValue 'E5' '  +0.56%  '
Set-TextValue 'D6' '145.75'
Set-TextValue 'E6' '  +3.65%  '
Set-TextValue 'E7' '  +0.09%  '
Set-TextValue 'D8' '0.573'
Set-TextValue 'E8' '  +0.90%  '
Set-TextValue 'D9' '2.669.31'
Set-TextValue 'E9' '  +2.24%  '
Set-TextValue 'D10' '6.66'
Set-TextValue 'E10' '  +2.97%  '
Set-TextValue 'E11' '  +0.66%  '
Set-TextValue 'D12' '0.338'
Set-TextValue 'E12' '  +0.43%  '
Set-TextValue 'E13' '  -0.73%  '
Set-TextValue 'D14' '3.128.71'
Set-TextValue 'E14' '  +2.30%  '
Set-TextValue 'D15' '59.666.25'
Set-TextValue 'E15' '  +0.66%  '
Set-TextValue 'D16' '21.21'
Set-TextValue 'E16' '  +3.15%  '
Set-TextValue 'D17' '2.714.60'
Set-TextValue 'E17' '  +5.19%  '
Set-TextValue 'E18' '  +1.21%  '
Set-TextValue 'D19' '344.24'
Set-TextValue 'E19' '  -0.51%  '
Set-TextValue 'D20' '4.42'
Set-TextValue 'E20' '  +1.88%  '
Set-TextValue 'D21' '10.43'
Set-TextValue 'E21' '  +2.99%  '
Set-TextValue 'D22' '6.33'
Set-TextValue 'E22' '  -0.57%  '
Set-TextValue 'D23' '0.998'
Set-TextValue 'E23' '  -0.12%  '
Set-TextValue 'D24' '66.79'
Set-TextValue 'E24' '  -0.73%  '
Set-TextValue 'E25' '  +2.32%  '
Set-TextValue 'E26' '  -1.30%  '
Set-TextValue 'E27' '  +0.10%  '
Set-TextValue 'E28' '  +1.03%  '
Set-TextValue 'D29' '0.0₃0758'
Set-TextValue 'E29' '  +2.58%  '
Set-TextValue 'D30' '0.998'
Set-TextValue 'E30' '  -0.04%  '
Set-TextValue 'E31' '  +1.75%  '
Set-TextValue 'D32' '5.86'
Set-TextValue 'E32' '  +0.44%  '
Set-TextValue 'D33' '19.02'
Set-TextValue 'E33' '  +0.82%  '
Set-TextValue 'D34' '150.23'
Set-TextValue 'E34' '  +0.60%  '
Set-TextValue 'E35' '  +0.91%  '
Set-TextValue 'E36' '  +2.20%  '
Set-TextValue 'B37' 'SuiNetwork'
Set-TextValue 'C37' 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue 'D37' '0.841'
Set-TextValue 'E37' '  -0.84%  '
Set-TextValue 'E38' '  -0.35%  '
Set-TextValue 'B39' 'Fetch.AI'
Set-TextValue 'C39' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D39' '0.844'
Set-TextValue 'E39' '  +0.72%  '
Set-TextValue 'D40' '293.51'
Set-TextValue 'E40' '  +5.83%  '
Set-TextValue 'D41' '3.60'
Set-TextValue 'E41' '  +1.46%  '
Set-TextValue 'E42' '  +0.12%  '
Set-TextValue 'E43' '  +1.44%  '
Set-TextValue 'D44' '19.54'
Set-TextValue 'E44' '  +4.86%  '
Set-TextValue 'D45' '0.0542'
Set-TextValue 'E45' '  +3.65%  '
Set-TextValue 'D46' '10.72'
Set-TextValue 'E46' '  -0.48%  '
Set-TextValue 'E47' '  -1.39%  '
Set-TextValue 'D48' '1.981.82'
Set-TextValue 'E48' '  +1.74%  '
Set-TextValue 'E49' '  +1.99%  '
Set-TextValue 'D50' '4.59'
Set-TextValue 'E50' '  +1.43%  '
Set-TextValue 'D51' '18.42'
Set-TextValue 'E51' '  -0.20%  '
